# Update "想去人数" (column F) counts across all sheets to match the
# regenerated gh-pages data snapshot (commit 456a3b4).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 558
$ws.Range("F3").Value = 256
$ws.Range("F4").Value = 584
$ws.Range("F5").Value = 1390
$ws.Range("F6").Value = 699
$ws.Range("F8").Value = 42
$ws.Range("F9").Value = 162
$ws.Range("F10").Value = 436
$ws.Range("F11").Value = 6441
$ws.Range("F13").Value = 30
$ws.Range("F14").Value = 1899
$ws.Range("F15").Value = 4782
$ws.Range("F18").Value = 5600
$ws.Range("F19").Value = 7525
$ws.Range("F20").Value = 150
$ws.Range("F22").Value = 768
$ws.Range("F23").Value = 4070
$ws.Range("F24").Value = 579
$ws.Range("F25").Value = 22
$ws.Range("F28").Value = 146
$ws.Range("F29").Value = 1073
$ws.Range("F30").Value = 1508
$ws.Range("F31").Value = 580
$ws.Range("F32").Value = 711
$ws.Range("F33").Value = 1710
$ws.Range("F34").Value = 248
$ws.Range("F35").Value = 1950
$ws.Range("F36").Value = 241
$ws.Range("F37").Value = 48
$ws.Range("F38").Value = 1267
$ws.Range("F39").Value = 1335
$ws.Range("F40").Value = 706
$ws.Range("F41").Value = 328
$ws.Range("F42").Value = 1929
$ws.Range("F43").Value = 3740
$ws.Range("F44").Value = 163
$ws.Range("F45").Value = 354
$ws.Range("F46").Value = 455
$ws.Range("F47").Value = 28
$ws.Range("F49").Value = 3972

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F3").Value = 1290
$ws.Range("F11").Value = 18
$ws.Range("F19").Value = 6
$ws.Range("F31").Value = 25

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 4528

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 4528
$ws.Range("F4").Value = 559
$ws.Range("F5").Value = 1290
$ws.Range("F8").Value = 256
$ws.Range("F9").Value = 584
$ws.Range("F11").Value = 1390
$ws.Range("F13").Value = 699
$ws.Range("F15").Value = 42
$ws.Range("F16").Value = 162
$ws.Range("F17").Value = 436
$ws.Range("F18").Value = 6441
$ws.Range("F20").Value = 4782
$ws.Range("F21").Value = 5600
$ws.Range("F22").Value = 5600
$ws.Range("F23").Value = 7525
$ws.Range("F25").Value = 768
$ws.Range("F26").Value = 4070
$ws.Range("F27").Value = 579
$ws.Range("F29").Value = 146
$ws.Range("F30").Value = 1073
$ws.Range("F31").Value = 1508
$ws.Range("F32").Value = 580
$ws.Range("F33").Value = 711
$ws.Range("F34").Value = 1710
$ws.Range("F35").Value = 248
$ws.Range("F36").Value = 1950
$ws.Range("F41").Value = 706
$ws.Range("F42").Value = 328
$ws.Range("F44").Value = 3740
$ws.Range("F45").Value = 25
$ws.Range("F46").Value = 163
$ws.Range("F47").Value = 354
$ws.Range("F50").Value = 3973
